$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('2026-02-12 17:35:02','237652194293','MABEKOU TALLA MARIELLE SANDRA _MINDEM SARL',429804),
    @('2026-02-12 14:47:02','237652386684','GUSTAVE KOUAM KUISSEU',12387),
    @('2026-02-12 17:35:53','237652927180','NDUMBE ROLINE JANONG ETS MOBILE FINANCIAL SERVICES MFS',206948),
    @('2026-02-12 16:07:17','237653650087','TCHOUALA THIERRY LTDLA_POLAS_BTQ_NDOGBONG ECOLE LAIC',170358),
    @('2026-02-12 16:14:57','237675396752','BENEDICTE CHANTAL MANTSANG',40772),
    @('2026-02-12 12:54:50','237675621967','NGO NTAMACK CECILE CECIC',6202),
    @('2026-02-12 15:30:41','237679548340','etienne ulrich awongo mbazoa',1786691),
    @('2026-02-12 16:06:04','237680144977','NDZEMBATEH NDZENAPU N ETS MOBILE FINANCIAL SERVICES MFS',47690),
    @('2026-02-12 18:28:26','237681679310','LA NEGRESSE SARL DIALLO AMADOU OURY',9125),
    @('2026-02-12 13:50:53','237682238745','RACHEL PRUDENCE JIKE KETCHA',3706014),
    @('2026-02-12 12:45:51','237682316602','GAELLE MAYEMELI ZAMPA',244667),
    @('2026-02-12 14:45:20','237682430965','ELSA CABRELLE MAKOUNGANG ETS MOBILE FINANCIAL SERVICES MFS',9296),
    @('2026-02-12 13:47:39','237683366333','ETS MOBILE FINANCIAL SERVICES MFS DJIAGUE JEANNETTE',37492),
    @('2026-02-12 13:09:36','237683379155','GUY LEONEL FOTSO SIMO',513256),
    @('2026-02-12 12:36:18','237683379207','christain alex mboo',157204),
    @('2026-02-12 15:39:10','237652049507','DONGMO GUESSOP ROVIELLE LAURA TOP MOBIL TELECOM',23056),
    @('2026-02-12 16:40:13','237671105715','DJOUMEGNI HORIELLE ETS MOBILE FINANCIAL SERVICES MFS',140383),
    @('2026-02-12 11:55:29','237674435926','FABRICE TCHOFFO',74695),
    @('2026-02-12 13:15:02','237674892352','FRANCOIS MITTERANA TAGUEGUIM',153550),
    @('2026-02-12 12:46:49','237674895389','LA NEGRESSE LTDLA_CBOX_R1_MABAH DASSI FANNY BLANCHE',16726),
    @('2026-02-12 17:01:25','237674908416','GNITEDEMO MEDJOUTSA CHARLOTTE SYLVIE ETS LE CONTENT',1912),
    @('2026-02-12 15:33:12','237675207106','germaine pougoue epse tchanang',1132334),
    @('2026-02-12 16:08:03','237675443028','LYSETTE LEONELLE  S NGUENA',663),
    @('2026-02-12 14:46:09','237676036869','KENFACK VIRGINIE GIODELE ETS MOBILE FINANCIAL SERVICES MFS',658203),
    @('2026-02-12 15:42:57','237677178575','GISELE CLAIRE ZAGBOUE DJIALA',74485),
    @('2026-02-12 17:04:33','237678183211','KAMTO NANGOP STEVE GAETAN ETS TCHATCHOUANG PAUL _ETP',27359),
    @('2026-02-12 16:18:25','237679038720','MUEINGNIEU MBIAWANG EPSE TCHATCHOUA ROSE ANNETTE Ets KENTO  Co',160874),
    @('2026-02-12 14:40:05','237679087694','Christ Brayant Penjeu Tchinkou',134023),
    @('2026-02-12 14:40:57','237679869809','XAVIEE ROSINE MEZAMO',15161),
    @('2026-02-12 16:23:52','237680435802','derice mboumela',835574),
    @('2026-02-12 12:21:39','237681114247','MINDEM SARL ymele voufoYMELE VOUFO VANESSA',68382),
    @('2026-02-12 07:45:38','237681602244','TSOMEJIO KENFACK NICAISE NESLIE ETS TCHATCHOUANG PAUL  ETP',509875),
    @('2026-02-12 14:10:39','237681606646','LA NEGRESSE SARL DONFACK PAULINE PELAGIE',100311),
    @('2026-02-12 15:16:00','237681655241','LA NEGRESSE LTDLA CBOX R1 NKUIDJEU KAMDOUM SYMPHORIEN',36416),
    @('2026-02-12 17:19:16','237681658403','LA NEGRESSE MISSOKE-UNIVERSITE',57370),
    @('2026-02-12 15:00:38','237681862876','TIDO GARLINE NOGRA-POLAS-BTQ-MAKEPE MISSOKE',24620),
    @('2026-02-10 10:49:58','237674243367','EMEGNI NGUEKAM DESTO WILFRIED CHIC MOBILE SARL',229),
    @('2026-02-12 17:34:32','237681656314','GATIEN TCHOUAYA MBOUCHIEKO',428865),
    @('2026-02-12 16:50:57','237650934256','DIDIER ROMUALD MBAKOP NYA',100918),
    @('2026-02-12 12:08:22','237652427111','CELESTINE CHANTAL MENDJOBOU EPSE NZIAKOU NJANJO',174411),
    @('2026-02-12 15:52:14','237672064755','KENFACK FRANC DUVIAL LA NEGRESSE SARL',1096209),
    @('2026-02-12 12:32:52','237673718583','Soppi Verole ETS MOBILE FINANCIAL SERVICES MFS',128142),
    @('2026-02-12 14:40:50','237675950748','ATANGANA HENRI EITEL ETS TCHATCHOUANG PAUL _ETP',1024033),
    @('2026-02-12 14:58:12','237683356603','MFS MATCHINDA SENDRINE',341),
    @('2026-02-12 16:25:52','237683395123','LA NEGRESSE SARL LTDLA-CBOX-R1-TAGNIN NICAISSE FLEURIE',58071),
    @('2026-02-12 17:27:21','237670174030','LA NEGRESSE LTDLA_CBOX_R1_FOHOM STEPHANE THIERRY',8250),
    @('2026-02-12 15:56:21','237681611433','Barry Diakariaou World T Plus',90790),
    @('2026-02-12 14:48:21','237652667691','MAPI TAKIAN LIVINE ETS MOBILE FINANCIAL SERVICES MFS',82481)
)

$startRow = 249

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
